$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3031.6897
$ws.Range("I38").Value = 325.58334
$ws.Range("J38").Value = 4941.8823
$ws.Range("K38").Value = 976.7500200000001
$ws.Range("L38").Value = 14825.6469
$ws.Range("M38").Value = -604.7500200000001
$ws.Range("N38").Value = -15569.6469
$ws.Range("H42").Value = 125000180
$ws.Range("I42").Value = 142857200
$ws.Range("K42").Value = 428571600
$ws.Range("M42").Value = -428571370
$ws.Range("H43").Value = 12960
$ws.Range("I43").Value = 9000
$ws.Range("J43").Value = 13950
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 13950
$ws.Range("M43").Value = -8931
$ws.Range("N43").Value = -14088
$ws.Range("H55").Value = 461.89474
$ws.Range("J55").Value = 511.53333
$ws.Range("L55").Value = 511.53333
$ws.Range("N55").Value = -939.53333
$ws.Range("H63").Value = 750018750
$ws.Range("J63").Value = 750018750
$ws.Range("L63").Value = 750018750
$ws.Range("N63").Value = -750019998
$ws.Range("H66").Value = 750018750
$ws.Range("J66").Value = 750018750
$ws.Range("L66").Value = 2250056250
$ws.Range("N66").Value = -2250062490
$ws.Range("H74").Value = 11288.235
$ws.Range("I74").Value = 11288.235
$ws.Range("K74").Value = 11288.235
$ws.Range("M74").Value = -10352.235
$ws.Range("H75").Value = 500033150
$ws.Range("J75").Value = 500033150
$ws.Range("L75").Value = 500033150
$ws.Range("N75").Value = -500035022
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H77").Value = 11288.235
$ws.Range("I77").Value = 11288.235
$ws.Range("K77").Value = 56441.175
$ws.Range("M77").Value = -51761.175
$ws.Range("H78").Value = 500033150
$ws.Range("J78").Value = 500033150
$ws.Range("L78").Value = 1500099450
$ws.Range("N78").Value = -1500108810
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H80").Value = 907.2941
$ws.Range("I80").Value = 627.7273
$ws.Range("J80").Value = 1419.8334
$ws.Range("K80").Value = 1883.1819
$ws.Range("L80").Value = 4259.5002
$ws.Range("M80").Value = -885.1819
$ws.Range("N80").Value = -6255.5002
$ws.Range("H83").Value = 907.2941
$ws.Range("I83").Value = 627.7273
$ws.Range("J83").Value = 1419.8334
$ws.Range("K83").Value = 5649.545700000001
$ws.Range("L83").Value = 12778.5006
$ws.Range("M83").Value = -657.5457000000006
$ws.Range("N83").Value = -22762.5006
$ws.Range("H86").Value = 3319.182
$ws.Range("I86").Value = 3813.875
$ws.Range("K86").Value = 3813.875
$ws.Range("M86").Value = -2690.875
$ws.Range("H87").Value = 600029950
$ws.Range("J87").Value = 600029950
$ws.Range("L87").Value = 600029950
$ws.Range("N87").Value = -600032446
$ws.Range("H88").Value = 6911
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H89").Value = 3319.182
$ws.Range("I89").Value = 3813.875
$ws.Range("K89").Value = 19069.375
$ws.Range("M89").Value = -13453.375
$ws.Range("H90").Value = 600029950
$ws.Range("J90").Value = 600029950
$ws.Range("L90").Value = 1800089850
$ws.Range("N90").Value = -1800102330
$ws.Range("H91").Value = 6911
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H98").Value = 8932.799999999999
$ws.Range("I98").Value = 10931.6875
$ws.Range("J98").Value = 937.25
$ws.Range("K98").Value = 10931.6875
$ws.Range("L98").Value = 937.25
$ws.Range("M98").Value = -9433.6875
$ws.Range("N98").Value = -3933.25
$ws.Range("H106").Value = 5350.9355
$ws.Range("J106").Value = 4523.8184
$ws.Range("L106").Value = 4523.8184
$ws.Range("N106").Value = -5785.8184
$ws.Range("H113").Value = 3873.7778
$ws.Range("I113").Value = 3644.1667
$ws.Range("J113").Value = 4333
$ws.Range("K113").Value = 3644.1667
$ws.Range("L113").Value = 4333
$ws.Range("M113").Value = -390.1667000000002
$ws.Range("N113").Value = -10841
$ws.Range("H122").Value = 8932.799999999999
$ws.Range("I122").Value = 10931.6875
$ws.Range("J122").Value = 937.25
$ws.Range("K122").Value = 32795.0625
$ws.Range("L122").Value = 2811.75
$ws.Range("M122").Value = -30345.0625
$ws.Range("N122").Value = -7711.75
$ws.Range("H132").Value = 1849.3611
$ws.Range("I132").Value = 1315.6
$ws.Range("K132").Value = 3946.8
$ws.Range("M132").Value = -1416.8
$ws.Range("H137").Value = 4021.2
$ws.Range("I137").Value = 2601.6667
$ws.Range("K137").Value = 7805.000100000001
$ws.Range("M137").Value = -5255.000100000001
$ws.Range("H138").Value = 3308.1267
$ws.Range("I138").Value = 1579.7142
$ws.Range("J138").Value = 4034.06
$ws.Range("K138").Value = 4739.142599999999
$ws.Range("L138").Value = 12102.18
$ws.Range("M138").Value = 400.8574000000008
$ws.Range("N138").Value = -22382.18
$ws.Range("H141").Value = 3716.75
$ws.Range("I141").Value = 3502.1875
$ws.Range("K141").Value = 10506.5625
$ws.Range("M141").Value = -5326.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3993.7
$ws.Range("I32").Value = 3453.5557
$ws.Range("K32").Value = 3453.5557
$ws.Range("M32").Value = -3166.5557
$ws.Range("H45").Value = 5249.75
$ws.Range("I45").Value = 4999.5
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 4999.5
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -4622.5
$ws.Range("N45").Value = -6254
$ws.Range("H61").Value = 5388.1763
$ws.Range("I61").Value = 5388.1763
$ws.Range("K61").Value = 5388.1763
$ws.Range("M61").Value = -5176.1763
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51748
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -158736
$ws.Range("H110").Value = 400
$ws.Range("I110").Value = 400
$ws.Range("K110").Value = 400
$ws.Range("M110").Value = 1645
$ws.Range("H122").Value = 4375.385
$ws.Range("I122").Value = 4685
$ws.Range("J122").Value = 3880
$ws.Range("K122").Value = 14055
$ws.Range("L122").Value = 11640
$ws.Range("M122").Value = -11605
$ws.Range("N122").Value = -16540
$ws.Range("H132").Value = 1630.5416
$ws.Range("I132").Value = 1687.8636
$ws.Range("K132").Value = 5063.5908
$ws.Range("M132").Value = -2533.5908
$ws.Range("H136").Value = 5388.1763
$ws.Range("I136").Value = 5388.1763
$ws.Range("K136").Value = 16164.5289
$ws.Range("M136").Value = -13614.5289

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 23500
$ws.Range("J49").Value = 23500
$ws.Range("L49").Value = 23500
$ws.Range("N49").Value = -23978
$ws.Range("H86").Value = 3065
$ws.Range("I86").Value = 2250.5
$ws.Range("K86").Value = 2250.5
$ws.Range("M86").Value = -1127.5
$ws.Range("H89").Value = 3065
$ws.Range("I89").Value = 2250.5
$ws.Range("K89").Value = 11252.5
$ws.Range("M89").Value = -5636.5
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H99").Value = 3009.5715
$ws.Range("I99").Value = 1445.1111
$ws.Range("J99").Value = 4182.9165
$ws.Range("K99").Value = 1445.1111
$ws.Range("L99").Value = 4182.9165
$ws.Range("M99").Value = 52.88889999999992
$ws.Range("N99").Value = -7178.9165
$ws.Range("H102").Value = 6414
$ws.Range("I102").Value = 6414
$ws.Range("K102").Value = 6414
$ws.Range("M102").Value = -3169
$ws.Range("H105").Value = 3411.8635
$ws.Range("I105").Value = 2863.3333
$ws.Range("J105").Value = 5880.25
$ws.Range("K105").Value = 2863.3333
$ws.Range("L105").Value = 5880.25
$ws.Range("M105").Value = -1116.3333
$ws.Range("N105").Value = -9374.25
$ws.Range("H107").Value = 4724.45
$ws.Range("I107").Value = 3793.4707
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 3793.4707
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = -1873.4707
$ws.Range("N107").Value = -13840
$ws.Range("H134").Value = 6487.25
$ws.Range("I134").Value = 6487.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 19461.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -16926.75
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1478.3
$ws.Range("I16").Value = 1222.875
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 1222.875
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -935.875
$ws.Range("N16").Value = -3074
$ws.Range("H31").Value = 3569.4194
$ws.Range("I31").Value = 1435.2
$ws.Range("J31").Value = 5570.25
$ws.Range("K31").Value = 1435.2
$ws.Range("L31").Value = 5570.25
$ws.Range("M31").Value = -1140.2
$ws.Range("N31").Value = -6160.25
$ws.Range("H34").Value = 3569.4194
$ws.Range("I34").Value = 1435.2
$ws.Range("J34").Value = 5570.25
$ws.Range("K34").Value = 1435.2
$ws.Range("L34").Value = 5570.25
$ws.Range("M34").Value = -1233.2
$ws.Range("N34").Value = -5974.25
$ws.Range("H36").Value = 2999.6667
$ws.Range("I36").Value = 2999.6667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2999.6667
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2611.6667
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 2999.6667
$ws.Range("I40").Value = 2999.6667
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2999.6667
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2839.6667
$ws.Range("N40").ClearContents()
$ws.Range("H58").Value = 2483.2222
$ws.Range("I58").Value = 2514.2144
$ws.Range("J58").Value = 2374.75
$ws.Range("K58").Value = 2514.2144
$ws.Range("L58").Value = 2374.75
$ws.Range("M58").Value = -2311.2144
$ws.Range("N58").Value = -2780.75
$ws.Range("H99").Value = 2891.6667
$ws.Range("I99").Value = 2891.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2891.6667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1393.6667
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 627.05554
$ws.Range("J107").Value = 733.125
$ws.Range("L107").Value = 733.125
$ws.Range("N107").Value = -4573.125
$ws.Range("H113").Value = 1478.3
$ws.Range("I113").Value = 1222.875
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1222.875
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 947.125
$ws.Range("N113").Value = -6840
$ws.Range("H118").Value = 75000
$ws.Range("J118").Value = 75000
$ws.Range("L118").Value = 75000
$ws.Range("N118").Value = -78314
$ws.Range("H122").Value = 1318.5264
$ws.Range("I122").Value = 1229
$ws.Range("K122").Value = 3687
$ws.Range("M122").Value = -1237
$ws.Range("H126").Value = 2891.6667
$ws.Range("I126").Value = 2891.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8675.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6205.000100000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1345.4615
$ws.Range("I132").Value = 1354.8857
$ws.Range("K132").Value = 4064.6571
$ws.Range("M132").Value = -1534.6571
$ws.Range("H134").Value = 2577.6843
$ws.Range("I134").Value = 2577.6843
$ws.Range("K134").Value = 7733.0529
$ws.Range("M134").Value = -5198.0529
$ws.Range("H136").Value = 2483.2222
$ws.Range("I136").Value = 2514.2144
$ws.Range("J136").Value = 2374.75
$ws.Range("K136").Value = 7542.6432
$ws.Range("L136").Value = 7124.25
$ws.Range("M136").Value = -4992.6432
$ws.Range("N136").Value = -12224.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 13468
$ws.Range("I8").Value = 13468
$ws.Range("K8").Value = 40404
$ws.Range("M8").Value = -40265
$ws.Range("H24").Value = 100
$ws.Range("I24").Value = 100
$ws.Range("K24").Value = 300
$ws.Range("M24").Value = -70
$ws.Range("H29").Value = 63.375
$ws.Range("J29").Value = 77.2
$ws.Range("L29").Value = 231.6
$ws.Range("N29").Value = -785.6
$ws.Range("H39").Value = 6025.2856
$ws.Range("J39").Value = 6999.6665
$ws.Range("L39").Value = 20998.9995
$ws.Range("N39").Value = -21586.9995
$ws.Range("H46").Value = 11212301
$ws.Range("I46").Value = 9956915
$ws.Range("K46").Value = 29870745
$ws.Range("M46").Value = -29870654
$ws.Range("H107").Value = 1425.5883
$ws.Range("I107").Value = 2390.6667
$ws.Range("J107").Value = 339.875
$ws.Range("K107").Value = 7172.000100000001
$ws.Range("L107").Value = 1019.625
$ws.Range("M107").Value = -5252.000100000001
$ws.Range("N107").Value = -4859.625
$ws.Range("H113").Value = 2354.1667
$ws.Range("J113").Value = 2929.7693
$ws.Range("L113").Value = 8789.3079
$ws.Range("N113").Value = -13129.3079
$ws.Range("H120").Value = 139394.9
$ws.Range("I120").Value = 217489.8
$ws.Range("K120").Value = 652469.3999999999
$ws.Range("M120").Value = -647631.3999999999
$ws.Range("H121").Value = 789.05554
$ws.Range("I121").Value = 529
$ws.Range("K121").Value = 1587
$ws.Range("M121").Value = -277
$ws.Range("H122").Value = 2082.5
$ws.Range("I122").Value = 698.6667
$ws.Range("J122").Value = 3466.3333
$ws.Range("K122").Value = 6288.0003
$ws.Range("L122").Value = 31196.9997
$ws.Range("M122").Value = -3838.0003
$ws.Range("N122").Value = -36096.9997
$ws.Range("H124").Value = 631.6667
$ws.Range("I124").Value = 631.6667
$ws.Range("K124").Value = 1895.0001
$ws.Range("M124").Value = 3014.9999
$ws.Range("H129").Value = 2396.353
$ws.Range("I129").Value = 267
$ws.Range("J129").Value = 6300.1665
$ws.Range("K129").Value = 801
$ws.Range("L129").Value = 18900.4995
$ws.Range("M129").Value = 4199
$ws.Range("N129").Value = -28900.4995
$ws.Range("H131").Value = 1710.2106
$ws.Range("I131").Value = 770
$ws.Range("J131").Value = 2093.2593
$ws.Range("K131").Value = 2310
$ws.Range("L131").Value = 6279.777900000001
$ws.Range("M131").Value = 2730
$ws.Range("N131").Value = -16359.7779
$ws.Range("H139").Value = 4898.2104
$ws.Range("I139").Value = 5758.25
$ws.Range("J139").Value = 4272.727
$ws.Range("K139").Value = 17274.75
$ws.Range("L139").Value = 12818.181
$ws.Range("M139").Value = -12134.75
$ws.Range("N139").Value = -23098.181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 27289.715
$ws.Range("J55").Value = 27500
$ws.Range("L55").Value = 27500
$ws.Range("N55").Value = -28154
$ws.Range("H80").Value = 3751.353
$ws.Range("J80").Value = 4670.381
$ws.Range("L80").Value = 4670.381
$ws.Range("N80").Value = -6666.381
$ws.Range("H83").Value = 3751.353
$ws.Range("J83").Value = 4670.381
$ws.Range("L83").Value = 23351.905
$ws.Range("N83").Value = -33335.905
$ws.Range("H113").Value = 6181.4546
$ws.Range("I113").Value = 6199.6
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 6199.6
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -4029.6
$ws.Range("N113").Value = -10340
$ws.Range("H122").Value = 3114.9666
$ws.Range("I122").Value = 2074.1428
$ws.Range("K122").Value = 6222.428400000001
$ws.Range("M122").Value = -3772.428400000001
$ws.Range("H132").Value = 3560.25
$ws.Range("I132").Value = 3078.8333
$ws.Range("K132").Value = 9236.499899999999
$ws.Range("M132").Value = -6706.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1539.1666
$ws.Range("I7").Value = 1213.6666
$ws.Range("K7").Value = 1213.6666
$ws.Range("M7").Value = -1101.6666
$ws.Range("H16").Value = 836.34784
$ws.Range("J16").Value = 2912.4
$ws.Range("L16").Value = 2912.4
$ws.Range("N16").Value = -3252.4
$ws.Range("H25").Value = 15000
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15460
$ws.Range("H40").Value = 4809.4346
$ws.Range("I40").Value = 3857.4443
$ws.Range("J40").Value = 5421.4287
$ws.Range("K40").Value = 3857.4443
$ws.Range("L40").Value = 5421.4287
$ws.Range("M40").Value = -3721.4443
$ws.Range("N40").Value = -5693.4287
$ws.Range("H46").Value = 1901.6666
$ws.Range("I46").Value = 1264.2
$ws.Range("J46").Value = 2698.5
$ws.Range("K46").Value = 1264.2
$ws.Range("L46").Value = 2698.5
$ws.Range("M46").Value = -1076.2
$ws.Range("N46").Value = -3074.5
$ws.Range("H61").Value = 3423.7878
$ws.Range("I61").Value = 2434.1738
$ws.Range("J61").Value = 5699.9
$ws.Range("K61").Value = 2434.1738
$ws.Range("L61").Value = 5699.9
$ws.Range("M61").Value = -2232.1738
$ws.Range("N61").Value = -6103.9
$ws.Range("H106").Value = 7832.6665
$ws.Range("J106").Value = 7832.6665
$ws.Range("L106").Value = 7832.6665
$ws.Range("N106").Value = -10356.6665
$ws.Range("H113").Value = 3423.7878
$ws.Range("I113").Value = 2434.1738
$ws.Range("J113").Value = 5699.9
$ws.Range("K113").Value = 2434.1738
$ws.Range("L113").Value = 5699.9
$ws.Range("M113").Value = -264.1738
$ws.Range("N113").Value = -10039.9
$ws.Range("H116").Value = 216192.8
$ws.Range("J116").Value = 216192.8
$ws.Range("L116").Value = 216192.8
$ws.Range("N116").Value = -225370.8
$ws.Range("H122").Value = 8063.273
$ws.Range("I122").Value = 6103.1333
$ws.Range("J122").Value = 12263.571
$ws.Range("K122").Value = 18309.3999
$ws.Range("L122").Value = 36790.713
$ws.Range("M122").Value = -15859.3999
$ws.Range("N122").Value = -41690.713
$ws.Range("H126").Value = 1539.1666
$ws.Range("I126").Value = 1213.6666
$ws.Range("K126").Value = 3640.9998
$ws.Range("M126").Value = -1170.9998
$ws.Range("H140").Value = 94750
$ws.Range("I140").Value = 94750
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 94750
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -89570
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 80000
$ws.Range("J27").Value = 80000
$ws.Range("L27").Value = 80000
$ws.Range("N27").Value = -80138
$ws.Range("H49").Value = 21625
$ws.Range("I49").Value = 10000
$ws.Range("K49").Value = 10000
$ws.Range("M49").Value = -9770
$ws.Range("H54").Value = 25000
$ws.Range("J54").Value = 25000
$ws.Range("L54").Value = 25000
$ws.Range("N54").Value = -26040
$ws.Range("H93").Value = 84966.664
$ws.Range("J93").Value = 89950
$ws.Range("L93").Value = 89950
$ws.Range("N93").Value = -94942
$ws.Range("H110").Value = 129124.375
$ws.Range("J110").Value = 129124.375
$ws.Range("L110").Value = 129124.375
$ws.Range("N110").Value = -137304.375
